$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.900.22'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '2.572.43'
$ws.Range('E3').Value = '  -2.45%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'516.56"
$ws.Range('E5').Value = '  -2.39%  '
$ws.Range('D6').Value = "'138.93"
$ws.Range('E6').Value = '  -4.00%  '
$ws.Range('D7').Value = "'0.997"
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = "'0.559"
$ws.Range('E8').Value = '  -1.87%  '
$ws.Range('D9').Value = '2.585.18'
$ws.Range('E9').Value = '  -2.72%  '
$ws.Range('D10').Value = "'6.44"
$ws.Range('E10').Value = '  -3.42%  '
$ws.Range('D11').Value = "'0.0993"
$ws.Range('E11').Value = '  -4.48%  '
$ws.Range('D12').Value = "'0.325"
$ws.Range('E12').Value = '  -3.75%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('D14').Value = '3.024.45'
$ws.Range('E14').Value = '  -2.62%  '
$ws.Range('D15').Value = '57.873.38'
$ws.Range('E15').Value = '  -2.28%  '
$ws.Range('D16').Value = "'20.04"
$ws.Range('E16').Value = '  -3.97%  '
$ws.Range('D17').Value = '2.566.08'
$ws.Range('E17').Value = '  -3.93%  '
$ws.Range('E18').Value = '  -3.99%  '
$ws.Range('D19').Value = "'333.20"
$ws.Range('E19').Value = '  -2.59%  '
$ws.Range('D20').Value = "'4.28"
$ws.Range('E20').Value = '  -4.18%  '
$ws.Range('D21').Value = "'10.06"
$ws.Range('E21').Value = '  -5.20%  '
$ws.Range('D22').Value = "'6.34"
$ws.Range('E22').Value = '  -0.67%  '
$ws.Range('D23').Value = "'0.999"
$ws.Range('E23').Value = '  +0.17%  '
$ws.Range('D24').Value = "'65.83"
$ws.Range('E24').Value = '  +0.46%  '
$ws.Range('E25').Value = '  -1.89%  '
$ws.Range('D26').Value = "'0.997"
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'0.398"
$ws.Range('E27').Value = '  -4.53%  '
$ws.Range('D28').Value = '2.688.26'
$ws.Range('E28').Value = '  -2.55%  '
$ws.Range('D29').Value = "'6.91"
$ws.Range('E29').Value = '  -4.34%  '
$ws.Range('D30').Value = "'0.998"
$ws.Range('E30').Value = '  -0.03%  '
$ws.Range('D31').Value = '0.0₃0712'
$ws.Range('E31').Value = '  -10.85%  '
$ws.Range('D32').Value = "'5.91"
$ws.Range('E32').Value = '  -8.07%  '
$ws.Range('E33').Value = '  -3.90%  '
$ws.Range('D34').Value = "'18.58"
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('D35').Value = "'149.17"
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').Value = "'3.89"
$ws.Range('E36').Value = '  -6.92%  '
$ws.Range('E37').Value = '  -7.26%  '
$ws.Range('D38').Value = "'36.11"
$ws.Range('E38').Value = '  -1.25%  '
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').Value = "'0.826"
$ws.Range('E39').Value = '  -6.02%  '
$ws.Range('B40').Value = 'Fetch.AI'
$ws.Range('C40').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D40').Value = "'0.829"
$ws.Range('E40').Value = '  -4.18%  '
$ws.Range('D41').Value = "'1.42"
$ws.Range('E41').Value = '  -4.59%  '
$ws.Range('D42').Value = "'3.50"
$ws.Range('E42').Value = '  -4.07%  '
$ws.Range('D43').Value = "'0.997"
$ws.Range('E43').Value = '  -0.07%  '
$ws.Range('D44').Value = "'272.89"
$ws.Range('E44').Value = '  +1.19%  '
$ws.Range('E45').Value = '  +0.19%  '
$ws.Range('D46').Value = "'0.587"
$ws.Range('E46').Value = '  -2.32%  '
$ws.Range('E47').Value = '  -3.79%  '
$ws.Range('D48').Value = "'0.0513"
$ws.Range('E48').Value = '  -4.52%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.968.26'
$ws.Range('E49').Value = '  -3.39%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').Value = "'18.30"
$ws.Range('E50').Value = '  -5.58%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = "'4.49"
$ws.Range('E51').Value = '  -4.79%  '
